$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "TZP"
$ws.Range("B12").Value = 9555.077980348269
$ws.Range("C12").Value = 9491.9783711322
$ws.Range("D12").Value = 0.00000000000003661285289874389
$ws.Range("E12").Value = "Model 2 (* Hospital)"
